$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on Price (D) and Volume(1h) (E) columns so that
# values like "1.014" or "0.9990" are stored verbatim as text, matching
# the source data (not converted to numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.653.64"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "1.859.15"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("D4").Value = "1.014"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "334.31"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").Value = "1.012"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").Value = "0.4612"
$ws.Range("E7").Value = "  -1.93%  "
$ws.Range("D8").Value = "0.3905"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("D9").Value = "46.56"
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("D10").Value = "0.07941"
$ws.Range("E10").Value = "  -1.77%  "
$ws.Range("D11").Value = "0.9990"
$ws.Range("E11").Value = "  -2.58%  "
$ws.Range("D12").Value = "21.61"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").Value = "1.866.47"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").Value = "5.946"
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("D15").Value = "7.206"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").Value = "88.23"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "0.06723"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").Value = "17.19"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("D21").Value = "1.011"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "27.662.86"
$ws.Range("E22").Value = "  -1.20%  "
$ws.Range("D23").Value = "5.430"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").Value = "2.312"
$ws.Range("E25").Value = "  -1.63%  "
$ws.Range("D26").Value = "2.085.93"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").Value = "159.55"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").Value = "19.64"
$ws.Range("E28").Value = "  -2.36%  "
$ws.Range("D29").Value = "2.127"
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("D30").Value = "5.455"
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("D31").Value = "121.74"
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("D32").Value = "0.9730"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").Value = "0.09384"
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("D34").Value = "3.630"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("D35").Value = "5.300"
$ws.Range("E35").Value = "  -1.70%  "
$ws.Range("D36").Value = "1.327"
$ws.Range("E36").Value = "  -6.04%  "
$ws.Range("D37").Value = "0.02234"
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("D38").Value = "0.06003"
$ws.Range("E38").Value = "  -2.46%  "
$ws.Range("D39").Value = "8.357"
$ws.Range("E39").Value = "  +3.27%  "
$ws.Range("D40").Value = "1.192"
$ws.Range("E40").Value = "  -2.67%  "
$ws.Range("D41").Value = "1.012"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").Value = "0.5931"
$ws.Range("E42").Value = "  -1.24%  "
$ws.Range("D43").Value = "0.1867"
$ws.Range("E43").Value = "  -1.54%  "
$ws.Range("D44").Value = "10.32"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "1.246"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("D46").Value = "0.5595"
$ws.Range("E46").Value = "  -2.27%  "
$ws.Range("D47").Value = "12.16"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("D48").Value = "1.913"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").Value = "0.06712"
$ws.Range("E49").Value = "  -3.41%  "
$ws.Range("D50").Value = "111.63"
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("D51").Value = "1.051"
$ws.Range("E51").Value = "  -2.17%  "
